$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (existing rows 18-24 shift down to 19-25).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly data point.
$ws.Cells.Item(18,1).Value  = 8
$ws.Cells.Item(18,2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(18,3).Value  = 'Coquimbo'
$ws.Cells.Item(18,4).Value  = 44809
$ws.Cells.Item(18,5).Value  = 4
$ws.Cells.Item(18,6).Value  = 100112026
$ws.Cells.Item(18,7).Value  = 'Haba'
$ws.Cells.Item(18,8).Value  = 'Sin especificar'
$ws.Cells.Item(18,9).Value  = 'Primera'
$ws.Cells.Item(18,10).Value = 520
$ws.Cells.Item(18,11).Value = 9500
$ws.Cells.Item(18,12).Value = 10000
$ws.Cells.Item(18,13).Value = 9750
$ws.Cells.Item(18,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18,15).Value = 'Provincia del Elquí'
$ws.Cells.Item(18,16).Value = 390
$ws.Cells.Item(18,17).Value = 25
$ws.Cells.Item(18,18).Value = 'Hortaliza'
